$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1977.2  # H92: 3952.6956 -> 1977.2
$ws.Cells.Item(92, 9).Value = 1784  # I92: 1832 -> 1784
$ws.Cells.Item(92, 10).Value = 2750  # J92: 8800 -> 2750
$ws.Cells.Item(92, 11).Value = 1784  # K92: 1832 -> 1784
$ws.Cells.Item(92, 12).Value = 2750  # L92: 8800 -> 2750
$ws.Cells.Item(92, 13).Value = -536  # M92: -584 -> -536
$ws.Cells.Item(92, 14).Value = -5246  # N92: -11296 -> -5246

$ws.Cells.Item(112, 8).Value = 1108.4807  # H112: 1165.3478 -> 1108.4807
$ws.Cells.Item(112, 9).Value = 465  # I112: 480 -> 465
$ws.Cells.Item(112, 10).Value = 1192.4131  # J112: 1248.9269 -> 1192.4131
$ws.Cells.Item(112, 11).Value = 1395  # K112: 1440 -> 1395
$ws.Cells.Item(112, 12).Value = 3577.2393  # L112: 3746.7807 -> 3577.2393
$ws.Cells.Item(112, 13).Value = -287  # M112: -332 -> -287
$ws.Cells.Item(112, 14).Value = -5793.2393  # N112: -5962.780699999999 -> -5793.2393

$ws.Cells.Item(138, 8).Value = 2799.15  # H138: 2920.158 -> 2799.15
$ws.Cells.Item(138, 9).Value = 2123.5833  # I138: 2271.182 -> 2123.5833
$ws.Cells.Item(138, 11).Value = 6370.749899999999  # K138: 6813.545999999999 -> 6370.749899999999
$ws.Cells.Item(138, 13).Value = -1230.749899999999  # M138: -1673.545999999999 -> -1230.749899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 0  # H25: 305.33334 -> 0
$ws.Cells.Item(25, 9).Value = 0  # I25: 208 -> 0
$ws.Cells.Item(25, 10).Value = 0  # J25: 500 -> 0
$ws.Cells.Item(25, 11).Value = 0  # K25: 208 -> 0
$ws.Cells.Item(25, 12).Value = 0  # L25: 500 -> 0
$ws.Cells.Item(25, 13).ClearContents()  # M25: 194 -> (removed)
$ws.Cells.Item(25, 14).ClearContents()  # N25: -1304 -> (removed)

$ws.Cells.Item(30, 8).Value = 2002.5  # H30: 1866.6666 -> 2002.5
$ws.Cells.Item(30, 9).Value = 0  # I30: 1866.6666 -> 0
$ws.Cells.Item(30, 10).Value = 2002.5  # J30: 0 -> 2002.5
$ws.Cells.Item(30, 11).Value = 0  # K30: 1866.6666 -> 0
$ws.Cells.Item(30, 12).Value = 2002.5  # L30: 0 -> 2002.5
$ws.Cells.Item(30, 13).ClearContents()  # M30: -1716.6666 -> (removed)
$ws.Cells.Item(30, 14).Value = -2302.5  # N30: None -> -2302.5

$ws.Cells.Item(61, 8).Value = 15626827  # H61: 16130940 -> 15626827
$ws.Cells.Item(61, 9).Value = 16668488  # I61: 17243246 -> 16668488
$ws.Cells.Item(61, 10).Value = 1907  # J61: 2507 -> 1907
$ws.Cells.Item(61, 11).Value = 16668488  # K61: 17243246 -> 16668488
$ws.Cells.Item(61, 12).Value = 1907  # L61: 2507 -> 1907
$ws.Cells.Item(61, 13).Value = -16668276  # M61: -17243034 -> -16668276
$ws.Cells.Item(61, 14).Value = -2331  # N61: -2931 -> -2331

$ws.Cells.Item(136, 8).Value = 15626827  # H136: 16130940 -> 15626827
$ws.Cells.Item(136, 9).Value = 16668488  # I136: 17243246 -> 16668488
$ws.Cells.Item(136, 10).Value = 1907  # J136: 2507 -> 1907
$ws.Cells.Item(136, 11).Value = 50005464  # K136: 51729738 -> 50005464
$ws.Cells.Item(136, 12).Value = 5721  # L136: 7521 -> 5721
$ws.Cells.Item(136, 13).Value = -50002914  # M136: -51727188 -> -50002914
$ws.Cells.Item(136, 14).Value = -10821  # N136: -12621 -> -10821

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1433.8334  # H92: 1503.2142 -> 1433.8334
$ws.Cells.Item(92, 9).Value = 1000  # I92: 925 -> 1000
$ws.Cells.Item(92, 10).Value = 1520.6  # J92: 1936.875 -> 1520.6
$ws.Cells.Item(92, 11).Value = 3000  # K92: 2775 -> 3000
$ws.Cells.Item(92, 12).Value = 4561.799999999999  # L92: 5810.625 -> 4561.799999999999
$ws.Cells.Item(92, 13).Value = -1752  # M92: -1527 -> -1752
$ws.Cells.Item(92, 14).Value = -7057.799999999999  # N92: -8306.625 -> -7057.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 10000  # H23: 2000 -> 10000
$ws.Cells.Item(23, 9).Value = 0  # I23: 2000 -> 0
$ws.Cells.Item(23, 10).Value = 10000  # J23: 0 -> 10000
$ws.Cells.Item(23, 11).Value = 0  # K23: 2000 -> 0
$ws.Cells.Item(23, 12).Value = 10000  # L23: 0 -> 10000
$ws.Cells.Item(23, 13).ClearContents()  # M23: -1777 -> (removed)
$ws.Cells.Item(23, 14).Value = -10446  # N23: None -> -10446

$ws.Cells.Item(34, 8).Value = 39800  # H34: 36200 -> 39800
$ws.Cells.Item(34, 10).Value = 39800  # J34: 36200 -> 39800
$ws.Cells.Item(34, 12).Value = 39800  # L34: 36200 -> 39800
$ws.Cells.Item(34, 14).Value = -40336  # N34: -36736 -> -40336

$ws.Cells.Item(64, 8).Value = 20000  # H64: 15000 -> 20000
$ws.Cells.Item(64, 10).Value = 20000  # J64: 15000 -> 20000
$ws.Cells.Item(64, 12).Value = 20000  # L64: 15000 -> 20000
$ws.Cells.Item(64, 14).Value = -20496  # N64: -15496 -> -20496

$ws.Cells.Item(67, 8).Value = 20000  # H67: 15000 -> 20000
$ws.Cells.Item(67, 10).Value = 20000  # J67: 15000 -> 20000
$ws.Cells.Item(67, 12).Value = 20000  # L67: 15000 -> 20000
$ws.Cells.Item(67, 14).Value = -21716  # N67: -16716 -> -21716

$ws.Cells.Item(70, 8).Value = 20656.076  # H70: 11205.629 -> 20656.076
$ws.Cells.Item(70, 9).Value = 28307.12  # I70: 13581.272 -> 28307.12
$ws.Cells.Item(70, 10).Value = 6993.5  # J70: 7185.3076 -> 6993.5
$ws.Cells.Item(70, 11).Value = 28307.12  # K70: 13581.272 -> 28307.12
$ws.Cells.Item(70, 12).Value = 6993.5  # L70: 7185.3076 -> 6993.5
$ws.Cells.Item(70, 13).Value = -28037.12  # M70: -13311.272 -> -28037.12
$ws.Cells.Item(70, 14).Value = -7533.5  # N70: -7725.3076 -> -7533.5

$ws.Cells.Item(73, 8).Value = 20656.076  # H73: 11205.629 -> 20656.076
$ws.Cells.Item(73, 9).Value = 28307.12  # I73: 13581.272 -> 28307.12
$ws.Cells.Item(73, 10).Value = 6993.5  # J73: 7185.3076 -> 6993.5
$ws.Cells.Item(73, 11).Value = 28307.12  # K73: 13581.272 -> 28307.12
$ws.Cells.Item(73, 12).Value = 6993.5  # L73: 7185.3076 -> 6993.5
$ws.Cells.Item(73, 13).Value = -27371.12  # M73: -12645.272 -> -27371.12
$ws.Cells.Item(73, 14).Value = -8865.5  # N73: -9057.3076 -> -8865.5

$ws.Cells.Item(76, 8).Value = 39800  # H76: 36200 -> 39800
$ws.Cells.Item(76, 10).Value = 39800  # J76: 36200 -> 39800
$ws.Cells.Item(76, 12).Value = 39800  # L76: 36200 -> 39800
$ws.Cells.Item(76, 14).Value = -40430  # N76: -36830 -> -40430

$ws.Cells.Item(79, 8).Value = 39800  # H79: 36200 -> 39800
$ws.Cells.Item(79, 10).Value = 39800  # J79: 36200 -> 39800
$ws.Cells.Item(79, 12).Value = 39800  # L79: 36200 -> 39800
$ws.Cells.Item(79, 14).Value = -41984  # N79: -38384 -> -41984

$ws.Cells.Item(113, 8).Value = 101472.7  # H113: 144380.86 -> 101472.7
$ws.Cells.Item(113, 9).Value = 200965.4  # I113: 500380.5 -> 200965.4
$ws.Cells.Item(113, 10).Value = 1980  # J113: 1981 -> 1980
$ws.Cells.Item(113, 11).Value = 200965.4  # K113: 500380.5 -> 200965.4
$ws.Cells.Item(113, 12).Value = 1980  # L113: 1981 -> 1980
$ws.Cells.Item(113, 13).Value = -198795.4  # M113: -498210.5 -> -198795.4
$ws.Cells.Item(113, 14).Value = -6320  # N113: -6321 -> -6320

$ws.Cells.Item(126, 8).Value = 5452.7646  # H126: 4692.2915 -> 5452.7646
$ws.Cells.Item(126, 9).Value = 3950  # I126: 3202.25 -> 3950
$ws.Cells.Item(126, 10).Value = 5653.1333  # J126: 5437.3125 -> 5653.1333
$ws.Cells.Item(126, 11).Value = 11850  # K126: 9606.75 -> 11850
$ws.Cells.Item(126, 12).Value = 16959.3999  # L126: 16311.9375 -> 16959.3999
$ws.Cells.Item(126, 13).Value = -9380  # M126: -7136.75 -> -9380
$ws.Cells.Item(126, 14).Value = -21899.3999  # N126: -21251.9375 -> -21899.3999

$ws.Cells.Item(132, 8).Value = 4951.875  # H132: 5302.1787 -> 4951.875
$ws.Cells.Item(132, 9).Value = 3946  # I132: 4348.294 -> 3946
$ws.Cells.Item(132, 10).Value = 6628.3335  # J132: 6776.364 -> 6628.3335
$ws.Cells.Item(132, 11).Value = 11838  # K132: 13044.882 -> 11838
$ws.Cells.Item(132, 12).Value = 19885.0005  # L132: 20329.092 -> 19885.0005
$ws.Cells.Item(132, 13).Value = -9308  # M132: -10514.882 -> -9308
$ws.Cells.Item(132, 14).Value = -24945.0005  # N132: -25389.092 -> -24945.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4833.2163  # H7: 5461.2583 -> 4833.2163
$ws.Cells.Item(7, 9).Value = 4939.4375  # I7: 6509.1816 -> 4939.4375
$ws.Cells.Item(7, 10).Value = 4752.2856  # J7: 4884.9 -> 4752.2856
$ws.Cells.Item(7, 11).Value = 4939.4375  # K7: 6509.1816 -> 4939.4375
$ws.Cells.Item(7, 12).Value = 4752.2856  # L7: 4884.9 -> 4752.2856
$ws.Cells.Item(7, 13).Value = -4827.4375  # M7: -6397.1816 -> -4827.4375
$ws.Cells.Item(7, 14).Value = -4976.2856  # N7: -5108.9 -> -4976.2856

$ws.Cells.Item(11, 8).Value = 0  # H11: 4500 -> 0
$ws.Cells.Item(11, 10).Value = 0  # J11: 4500 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 4500 -> 0
$ws.Cells.Item(11, 14).ClearContents()  # N11: -4780 -> (removed)

$ws.Cells.Item(16, 8).Value = 3848.875  # H16: 3859.1 -> 3848.875
$ws.Cells.Item(16, 9).Value = 3848.875  # I16: 3555.8572 -> 3848.875
$ws.Cells.Item(16, 10).Value = 0  # J16: 4566.6665 -> 0
$ws.Cells.Item(16, 11).Value = 3848.875  # K16: 3555.8572 -> 3848.875
$ws.Cells.Item(16, 12).Value = 0  # L16: 4566.6665 -> 0
$ws.Cells.Item(16, 13).Value = -3678.875  # M16: -3385.8572 -> -3678.875
$ws.Cells.Item(16, 14).ClearContents()  # N16: -4906.6665 -> (removed)

$ws.Cells.Item(68, 8).Value = 2521.8333  # H68: 2366.6667 -> 2521.8333
$ws.Cells.Item(68, 9).Value = 2185  # I68: 2350 -> 2185
$ws.Cells.Item(68, 10).Value = 2858.6667  # J68: 2400 -> 2858.6667
$ws.Cells.Item(68, 11).Value = 2185  # K68: 2350 -> 2185
$ws.Cells.Item(68, 12).Value = 2858.6667  # L68: 2400 -> 2858.6667
$ws.Cells.Item(68, 13).Value = -1436  # M68: -1601 -> -1436
$ws.Cells.Item(68, 14).Value = -4356.6667  # N68: -3898 -> -4356.6667

$ws.Cells.Item(71, 8).Value = 2521.8333  # H71: 2366.6667 -> 2521.8333
$ws.Cells.Item(71, 9).Value = 2185  # I71: 2350 -> 2185
$ws.Cells.Item(71, 10).Value = 2858.6667  # J71: 2400 -> 2858.6667
$ws.Cells.Item(71, 11).Value = 10925  # K71: 11750 -> 10925
$ws.Cells.Item(71, 12).Value = 14293.3335  # L71: 12000 -> 14293.3335
$ws.Cells.Item(71, 13).Value = -7181  # M71: -8006 -> -7181
$ws.Cells.Item(71, 14).Value = -21781.3335  # N71: -19488 -> -21781.3335

$ws.Cells.Item(122, 8).Value = 5068.56  # H122: 4507.4136 -> 5068.56
$ws.Cells.Item(122, 9).Value = 6561  # I122: 4972.2144 -> 6561
$ws.Cells.Item(122, 11).Value = 19683  # K122: 14916.6432 -> 19683
$ws.Cells.Item(122, 13).Value = -17233  # M122: -12466.6432 -> -17233

$ws.Cells.Item(126, 8).Value = 4833.2163  # H126: 5461.2583 -> 4833.2163
$ws.Cells.Item(126, 9).Value = 4939.4375  # I126: 6509.1816 -> 4939.4375
$ws.Cells.Item(126, 10).Value = 4752.2856  # J126: 4884.9 -> 4752.2856
$ws.Cells.Item(126, 11).Value = 14818.3125  # K126: 19527.5448 -> 14818.3125
$ws.Cells.Item(126, 12).Value = 14256.8568  # L126: 14654.7 -> 14256.8568
$ws.Cells.Item(126, 13).Value = -12348.3125  # M126: -17057.5448 -> -12348.3125
$ws.Cells.Item(126, 14).Value = -19196.8568  # N126: -19594.7 -> -19196.8568

$ws.Cells.Item(136, 8).Value = 62518988  # H136: 41680544 -> 62518988
$ws.Cells.Item(136, 9).Value = 100002780  # I136: 62502940 -> 100002780
$ws.Cells.Item(136, 10).Value = 46001.668  # J136: 35751.25 -> 46001.668
$ws.Cells.Item(136, 11).Value = 300008340  # K136: 187508820 -> 300008340
$ws.Cells.Item(136, 12).Value = 138005.004  # L136: 107253.75 -> 138005.004
$ws.Cells.Item(136, 13).Value = -300005790  # M136: -187506270 -> -300005790
$ws.Cells.Item(136, 14).Value = -143105.004  # N136: -112353.75 -> -143105.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0  # H24: 11000 -> 0
$ws.Cells.Item(24, 10).Value = 0  # J24: 11000 -> 0
$ws.Cells.Item(24, 12).Value = 0  # L24: 11000 -> 0
$ws.Cells.Item(24, 14).ClearContents()  # N24: -11460 -> (removed)
